$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("68:68").Insert()

$ws.Range("A68").Value = 10
$ws.Range("B68").Value = "Vega Modelo de Temuco"
$ws.Range("C68").Value = "La Araucanía"
$ws.Range("D68").Value = 45225
$ws.Range("E68").Value = 9
$ws.Range("F68").Value = 100112010
$ws.Range("G68").Value = "Achicoria"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 150
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = 10000
$ws.Range("N68").Value = "$/caja 18 unidades"
$ws.Range("O68").Value = "Región Metropolitana"
$ws.Range("P68").Value = 556
$ws.Range("Q68").Value = 18
$ws.Range("R68").Value = "Hortaliza"
